# xlsx mit constraints ergaenzt
# Adds two new "LG_GeoAssets_V2" constraint rows (NatRel) plus a new
# "GeolAssets_V2" association constraint row to the Constraints sheet,
# and updates the active sheet/selection bookkeeping accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constraints")

# Insert one new row before row 13; this pushes the existing rows 13-38
# down to 14-39 and keeps their content/styles/heights intact.
$ws.Rows.Item(13).Insert()

# --- Fill row 18 (new "NatRel" constraint #17) -----------------------
# Pick up the same cell styling used by the other "LG_GeoAssets_V2 /
# class / LGAssetItem" rows (copy format only from row 14, A:G).
$ws.Range("A14:G14").Copy()
$ws.Range("A18:G18").PasteSpecial(-4122)
$ws.Range("A18").Value2 = 17
$ws.Range("B18").Value2 = "LG_GeoAssets_V2"
$ws.Range("C18").Value2 = "class"
$ws.Range("D18").Value2 = "LGAssetItem"
$ws.Range("E18").Value2 = "Wenn das Asset nicht von nationaler Relevanz ist, muss dein Typ definiert werden"
$ws.Range("F18").Value2 = "SET CONSTRAINT WHERE NOT(IsNatRel): NOT (DEFINED (TypeNatRel));"
$ws.Range("G18").Value2 = "active"
$ws.Rows.Item(18).RowHeight = 25

# --- Fill row 19 (new "NatRel" constraint #18) -----------------------
# Same styling, but this row has no "status" (G) column entry at all.
$ws.Range("A14:F14").Copy()
$ws.Range("A19:F19").PasteSpecial(-4122)
$ws.Range("A19").Value2 = 18
$ws.Range("B19").Value2 = "LG_GeoAssets_V2"
$ws.Range("C19").Value2 = "class"
$ws.Range("D19").Value2 = "LGAssetItem"
$ws.Range("E19").Value2 = "Wenn das Asset nicht von nationaler Relevanz ist, muss dein Typ definiert werden"
$ws.Range("F19").Value2 = "SET CONSTRAINT WHERE (IsNatRel): DEFINED(TypeNatRel);"
$ws.Rows.Item(19).RowHeight = 25

# --- Fill the newly inserted row 13 (new "AssetItemMain_AssetItemPart"
#     association constraint #12) ------------------------------------
# This row uses the "association" header styling (same as rows 2-12).
$ws.Range("A2:F2").Copy()
$ws.Range("A13:F13").PasteSpecial(-4122)
$ws.Range("A13").Value2 = 12
$ws.Range("B13").Value2 = "GeolAssets_V2"
$ws.Range("C13").Value2 = "association"
$ws.Range("D13").Value2 = "AssetItemMain_AssetItemPart"
$ws.Range("F13").Value2 = "SET CONSTRAINT WHERE AssetItemPart->IsExtract: DEFINED (AssetItemMain);"
$ws.Range("E13").Value2 = "Jedes AssetPart, dass ein Extract ist, muss ein Main besitzen"
$ws.Rows.Item(13).RowHeight = 25

# Remove the leftover G13/H13 cells the row-insert carried over from
# row 12 (copied style+blank value); the new row has neither.
$ws.Range("G13").Clear()
$ws.Range("H13").Clear()

# --- Update active-sheet / selection bookkeeping ----------------------
$ws.Activate()
$ws.Range("D12").Select()
